$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Rocky Fork Middle School (sheet7.xml)
#   - rows 13-15: "### (if needed)" strings replaced by plain numbers
#   - row 16 added: "Fake 1" / 30
# ------------------------------------------------------------------
$wsRockyFork = $wb.Worksheets.Item("Rocky Fork Middle School")
$wsRockyFork.Range("A13").Value = 202
$wsRockyFork.Range("A14").Value = 205
$wsRockyFork.Range("A15").Value = 220
$wsRockyFork.Range("A16").Value = "Fake 1"
$wsRockyFork.Range("D16").Value = 30

# ------------------------------------------------------------------
# Siegel Middle School (sheet8.xml)
#   - rows 17-20 added: Fake 1..4 / 25, carrying row 16's formatting
# ------------------------------------------------------------------
$wsSiegel = $wb.Worksheets.Item("Siegel Middle School")

[void]$wsSiegel.Range("A16").Copy()
[void]$wsSiegel.Range("A17").PasteSpecial(-4122)
[void]$wsSiegel.Range("D16").Copy()
[void]$wsSiegel.Range("D17").PasteSpecial(-4122)
$wsSiegel.Range("A17").Value = "Fake 1"
$wsSiegel.Range("D17").Value = 25

[void]$wsSiegel.Range("A17").Copy()
[void]$wsSiegel.Range("A18").PasteSpecial(-4122)
[void]$wsSiegel.Range("D17").Copy()
[void]$wsSiegel.Range("D18").PasteSpecial(-4122)
$wsSiegel.Range("A18").Value = "Fake 2"
$wsSiegel.Range("D18").Value = 25

[void]$wsSiegel.Range("A17").Copy()
[void]$wsSiegel.Range("A19").PasteSpecial(-4122)
[void]$wsSiegel.Range("D17").Copy()
[void]$wsSiegel.Range("D19").PasteSpecial(-4122)
$wsSiegel.Range("A19").Value = "Fake 3"
$wsSiegel.Range("D19").Value = 25

[void]$wsSiegel.Range("A17").Copy()
[void]$wsSiegel.Range("A20").PasteSpecial(-4122)
[void]$wsSiegel.Range("D17").Copy()
[void]$wsSiegel.Range("D20").PasteSpecial(-4122)
$wsSiegel.Range("A20").Value = "Fake 4"
$wsSiegel.Range("D20").Value = 25

# ------------------------------------------------------------------
# Selection / active-tab bookkeeping.
# Siegel's selection moves to E24 but Siegel never becomes the active
# sheet, so update its selection first; Rocky Fork becomes the active
# sheet/tab (was Smyrna before), so select it last so it "wins".
# ------------------------------------------------------------------
[void]$wsSiegel.Range("E24").Select()
[void]$wsRockyFork.Range("D17").Select()

Write-Host "edit complete"
